$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand new row above row 234, shifting existing rows 234:284 down to 235:285.
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new record.
$ws.Range("A234").Value = 4
$ws.Range("B234").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C234").Value = "Los Lagos"
$ws.Range("D234").Value = 44798
$ws.Range("D234").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E234").Value = 10
$ws.Range("F234").Value = 100112032
$ws.Range("G234").Value = "Zapallo italiano"
$ws.Range("H234").Value = "Sin especificar"
$ws.Range("I234").Value = "Primera"
$ws.Range("J234").Value = 70
$ws.Range("K234").Value = 27000
$ws.Range("L234").Value = 27000
$ws.Range("M234").Value = 27000
$ws.Range("N234").Value = "`$/caja 50 unidades"
$ws.Range("O234").Value = "Región de Arica y Parinacota"
$ws.Range("P234").Value = 540
$ws.Range("Q234").Value = 50
$ws.Range("R234").Value = "Hortaliza"
